# Add the new "UNIQUE" constraint entry to the categories table (Sheet1!E3),
# matching the formatting used by the other populated cells in that column
# block, auto-fit the two newly-used columns (D, E) so their width reflects
# the content, and leave the selection where the author left it (E4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new constraint text - this also grows the shared-string table.
$ws.Range("E3").Value = "UNIQUE"

# Match the look of the neighbouring "value" cells (e.g. C3/C4 - Arial 10,
# black) instead of leaving the default Calibri formatting behind.
$ws.Range("E3").Font.Name = "Arial"
$ws.Range("E3").Font.Size = 10
$ws.Range("E3").Font.Color = 0

# Columns D and E now hold real content worth sizing to fit.
$ws.Columns.Item(4).AutoFit() | Out-Null
$ws.Columns.Item(5).AutoFit() | Out-Null

# Leave the cursor on E4, where the author ended up after the edit.
$ws.Range("E4").Select() | Out-Null
